$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.026.89"
$ws.Range("E2").Value = "  -3.95%  "
$ws.Range("D3").Value = "3.658.12"
$ws.Range("E3").Value = "  -4.38%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.57"
$ws.Range("E5").Value = "  -3.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.03"
$ws.Range("E6").Value = "  -6.55%  "
$ws.Range("D7").Value = "3.658.78"
$ws.Range("E7").Value = "  -4.29%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  -1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").Value = "  -5.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.13"
$ws.Range("E11").Value = "  -5.15%  "
$ws.Range("E12").Value = "  -5.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.65"
$ws.Range("E13").Value = "  -6.22%  "
$ws.Range("E14").Value = "  -6.75%  "
$ws.Range("D15").Value = "4.281.52"
$ws.Range("E15").Value = "  -3.91%  "
$ws.Range("D16").Value = "3.668.74"
$ws.Range("E16").Value = "  -3.88%  "
$ws.Range("D17").Value = "67.117.51"
$ws.Range("E17").Value = "  -3.95%  "
$ws.Range("E18").Value = "  -4.33%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.09"
$ws.Range("E19").Value = "  -6.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.75"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.55"
$ws.Range("E21").Value = "  -4.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.03"
$ws.Range("E22").Value = "  -7.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.713"
$ws.Range("E23").Value = "  -3.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.92"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("E25").Value = "  -7.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000138"
$ws.Range("E26").Value = "  -4.81%  "
$ws.Range("E27").Value = "  -5.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.92"
$ws.Range("E29").Value = "  -6.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.90"
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("E31").Value = "  -6.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.71"
$ws.Range("E32").Value = "  -4.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.61"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("D34").Value = "3.803.27"
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.106"
$ws.Range("E35").Value = "  -7.46%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.598.91"
$ws.Range("E36").Value = "  -4.11%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.979"
$ws.Range("E38").Value = "  -6.87%  "
$ws.Range("E39").Value = "  -6.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("E40").Value = "  -7.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.320"
$ws.Range("E41").Value = "  -6.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "439.85"
$ws.Range("E42").Value = "  -9.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.42"
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("E44").Value = "  -7.89%  "
$ws.Range("E45").Value = "  -8.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.29"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "141.50"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "39.51"
$ws.Range("E49").Value = "  -10.92%  "
$ws.Range("D50").Value = "2.745.73"
$ws.Range("E50").Value = "  -7.23%  "
$ws.Range("E51").Value = "  -5.38%  "
